# Scheduled market-data refresh for the Leve profit sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Updates currentAveragePrice / NQ / HQ price
# columns (H-L) and the derived LeveProfitNQ/HQ columns (M, N) with the
# latest pulled values. A handful of rows drop their M/N profit cell
# entirely when the corresponding price inputs are no longer available
# (e.g. NQ or HQ side reports 0), matching how the refresh script skips
# writing an undefined profit; a couple of rows gain a profit cell where
# one is now computable.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 231.25  # H2: 227.77777 -> 231.25
$ws.Cells.Item(2, 9).Value = 233.33333  # I2: 228.57143 -> 233.33333
$ws.Cells.Item(2, 11).Value = 233.33333  # K2: 228.57143 -> 233.33333
$ws.Cells.Item(2, 13).Value = -120.33333  # M2: -115.57143 -> -120.33333
$ws.Cells.Item(4, 8).Value = 235.5  # H4: 240.33333 -> 235.5
$ws.Cells.Item(4, 9).Value = 235.5  # I4: 240.33333 -> 235.5
$ws.Cells.Item(4, 11).Value = 235.5  # K4: 240.33333 -> 235.5
$ws.Cells.Item(4, 13).Value = -121.5  # M4: -126.33333 -> -121.5
$ws.Cells.Item(21, 8).Value = 0  # H21: 1000 -> 0
$ws.Cells.Item(21, 9).Value = 0  # I21: 1000 -> 0
$ws.Cells.Item(21, 11).Value = 0  # K21: 1000 -> 0
$ws.Cells.Item(21, 13).Value = $null  # M21: -532 -> None
$ws.Cells.Item(23, 8).Value = 0  # H23: 1000 -> 0
$ws.Cells.Item(23, 9).Value = 0  # I23: 1000 -> 0
$ws.Cells.Item(23, 11).Value = 0  # K23: 1000 -> 0
$ws.Cells.Item(23, 13).Value = $null  # M23: -766 -> None
$ws.Cells.Item(33, 8).Value = 466.7647  # H33: 506.8 -> 466.7647
$ws.Cells.Item(33, 9).Value = 466.7143  # I33: 497.30768 -> 466.7143
$ws.Cells.Item(33, 10).Value = 467  # J33: 568.5 -> 467
$ws.Cells.Item(33, 11).Value = 466.7143  # K33: 497.30768 -> 466.7143
$ws.Cells.Item(33, 12).Value = 467  # L33: 568.5 -> 467
$ws.Cells.Item(33, 13).Value = -237.7143  # M33: -268.30768 -> -237.7143
$ws.Cells.Item(33, 14).Value = -925  # N33: -1026.5 -> -925
$ws.Cells.Item(38, 8).Value = 134.22223  # H38: 436.16666 -> 134.22223
$ws.Cells.Item(38, 9).Value = 134.22223  # I38: 121.6 -> 134.22223
$ws.Cells.Item(38, 10).Value = 0  # J38: 2009 -> 0
$ws.Cells.Item(38, 11).Value = 402.66669  # K38: 364.8 -> 402.66669
$ws.Cells.Item(38, 12).Value = 0  # L38: 6027 -> 0
$ws.Cells.Item(38, 13).Value = -30.66669000000002  # M38: 7.200000000000045 -> -30.66669000000002
$ws.Cells.Item(38, 14).Value = $null  # N38: -6771 -> None
$ws.Cells.Item(39, 8).Value = 824.2857  # H39: 413.70587 -> 824.2857
$ws.Cells.Item(39, 9).Value = 354  # I39: 198.78572 -> 354
$ws.Cells.Item(39, 10).Value = 2000  # J39: 1416.6666 -> 2000
$ws.Cells.Item(39, 11).Value = 1062  # K39: 596.35716 -> 1062
$ws.Cells.Item(39, 12).Value = 6000  # L39: 4249.9998 -> 6000
$ws.Cells.Item(39, 13).Value = -766  # M39: -300.35716 -> -766
$ws.Cells.Item(39, 14).Value = -6592  # N39: -4841.9998 -> -6592
$ws.Cells.Item(55, 8).Value = 1243.8  # H55: 1244.2 -> 1243.8
$ws.Cells.Item(55, 9).Value = 419.42856  # I55: 423.33334 -> 419.42856
$ws.Cells.Item(55, 10).Value = 3167.3333  # J55: 2475.5 -> 3167.3333
$ws.Cells.Item(55, 11).Value = 419.42856  # K55: 423.33334 -> 419.42856
$ws.Cells.Item(55, 12).Value = 3167.3333  # L55: 2475.5 -> 3167.3333
$ws.Cells.Item(55, 13).Value = -205.42856  # M55: -209.33334 -> -205.42856
$ws.Cells.Item(55, 14).Value = -3595.3333  # N55: -2903.5 -> -3595.3333
$ws.Cells.Item(62, 8).Value = 7062.0586  # H62: 7390.467 -> 7062.0586
$ws.Cells.Item(62, 9).Value = 7796.769  # I62: 8055 -> 7796.769
$ws.Cells.Item(62, 10).Value = 4674.25  # J62: 4732.3335 -> 4674.25
$ws.Cells.Item(62, 11).Value = 7796.769  # K62: 8055 -> 7796.769
$ws.Cells.Item(62, 12).Value = 4674.25  # L62: 4732.3335 -> 4674.25
$ws.Cells.Item(62, 13).Value = -7172.769  # M62: -7431 -> -7172.769
$ws.Cells.Item(62, 14).Value = -5922.25  # N62: -5980.3335 -> -5922.25
$ws.Cells.Item(64, 8).Value = 8390.916999999999  # H64: 7762.0713 -> 8390.916999999999
$ws.Cells.Item(64, 9).Value = 3963.6667  # I64: 3973.8 -> 3963.6667
$ws.Cells.Item(64, 11).Value = 3963.6667  # K64: 3973.8 -> 3963.6667
$ws.Cells.Item(64, 13).Value = -3715.6667  # M64: -3725.8 -> -3715.6667
$ws.Cells.Item(65, 8).Value = 7062.0586  # H65: 7390.467 -> 7062.0586
$ws.Cells.Item(65, 9).Value = 7796.769  # I65: 8055 -> 7796.769
$ws.Cells.Item(65, 10).Value = 4674.25  # J65: 4732.3335 -> 4674.25
$ws.Cells.Item(65, 11).Value = 38983.845  # K65: 40275 -> 38983.845
$ws.Cells.Item(65, 12).Value = 23371.25  # L65: 23661.6675 -> 23371.25
$ws.Cells.Item(65, 13).Value = -35863.845  # M65: -37155 -> -35863.845
$ws.Cells.Item(65, 14).Value = -29611.25  # N65: -29901.6675 -> -29611.25
$ws.Cells.Item(67, 8).Value = 8390.916999999999  # H67: 7762.0713 -> 8390.916999999999
$ws.Cells.Item(67, 9).Value = 3963.6667  # I67: 3973.8 -> 3963.6667
$ws.Cells.Item(67, 11).Value = 3963.6667  # K67: 3973.8 -> 3963.6667
$ws.Cells.Item(67, 13).Value = -3105.6667  # M67: -3115.8 -> -3105.6667
$ws.Cells.Item(69, 8).Value = 7000  # H69: 6999.6665 -> 7000
$ws.Cells.Item(69, 9).Value = 4000  # I69: 0 -> 4000
$ws.Cells.Item(69, 10).Value = 10000  # J69: 6999.6665 -> 10000
$ws.Cells.Item(69, 11).Value = 12000  # K69: 0 -> 12000
$ws.Cells.Item(69, 12).Value = 30000  # L69: 20998.9995 -> 30000
$ws.Cells.Item(69, 13).Value = -11126  # M69: None -> -11126
$ws.Cells.Item(69, 14).Value = -31748  # N69: -22746.9995 -> -31748
$ws.Cells.Item(72, 8).Value = 7000  # H72: 6999.6665 -> 7000
$ws.Cells.Item(72, 9).Value = 4000  # I72: 0 -> 4000
$ws.Cells.Item(72, 10).Value = 10000  # J72: 6999.6665 -> 10000
$ws.Cells.Item(72, 11).Value = 36000  # K72: 0 -> 36000
$ws.Cells.Item(72, 12).Value = 90000  # L72: 62996.9985 -> 90000
$ws.Cells.Item(72, 13).Value = -31632  # M72: None -> -31632
$ws.Cells.Item(72, 14).Value = -98736  # N72: -71732.9985 -> -98736
$ws.Cells.Item(138, 8).Value = 2226.4194  # H138: 2226.7937 -> 2226.4194
$ws.Cells.Item(138, 10).Value = 3173.7036  # J138: 3140.7144 -> 3173.7036
$ws.Cells.Item(138, 12).Value = 9521.110799999999  # L138: 9422.143199999999 -> 9521.110799999999
$ws.Cells.Item(138, 14).Value = -19801.1108  # N138: -19702.1432 -> -19801.1108
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 1553.4  # H35: 1567 -> 1553.4
$ws.Cells.Item(35, 10).Value = 1499.75  # J35: 1500 -> 1499.75
$ws.Cells.Item(35, 12).Value = 1499.75  # L35: 1500 -> 1499.75
$ws.Cells.Item(35, 14).Value = -2311.75  # N35: -2312 -> -2311.75
$ws.Cells.Item(122, 8).Value = 3150.1538  # H122: 3395.3 -> 3150.1538
$ws.Cells.Item(122, 9).Value = 2989.8  # I122: 3271.2856 -> 2989.8
$ws.Cells.Item(122, 11).Value = 8969.400000000001  # K122: 9813.856800000001 -> 8969.400000000001
$ws.Cells.Item(122, 13).Value = -6519.400000000001  # M122: -7363.856800000001 -> -6519.400000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(27, 8).Value = 0  # H27: 145000 -> 0
$ws.Cells.Item(27, 10).Value = 0  # J27: 145000 -> 0
$ws.Cells.Item(27, 12).Value = 0  # L27: 145000 -> 0
$ws.Cells.Item(27, 14).Value = $null  # N27: -145384 -> None
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4480.2354  # H31: 4525.091 -> 4480.2354
$ws.Cells.Item(31, 10).Value = 5821.6665  # J31: 6023.2144 -> 5821.6665
$ws.Cells.Item(31, 12).Value = 5821.6665  # L31: 6023.2144 -> 5821.6665
$ws.Cells.Item(31, 14).Value = -6411.6665  # N31: -6613.2144 -> -6411.6665
$ws.Cells.Item(34, 8).Value = 4480.2354  # H34: 4525.091 -> 4480.2354
$ws.Cells.Item(34, 10).Value = 5821.6665  # J34: 6023.2144 -> 5821.6665
$ws.Cells.Item(34, 12).Value = 5821.6665  # L34: 6023.2144 -> 5821.6665
$ws.Cells.Item(34, 14).Value = -6225.6665  # N34: -6427.2144 -> -6225.6665
$ws.Cells.Item(97, 8).Value = 33406.168  # H97: 33421.332 -> 33406.168
$ws.Cells.Item(97, 10).Value = 34087.4  # J97: 34105.6 -> 34087.4
$ws.Cells.Item(97, 12).Value = 34087.4  # L97: 34105.6 -> 34087.4
$ws.Cells.Item(97, 14).Value = -36069.4  # N97: -36087.6 -> -36069.4
$ws.Cells.Item(139, 8).Value = 142876  # H139: 143176 -> 142876
$ws.Cells.Item(139, 9).Value = 0  # I139: 139000 -> 0
$ws.Cells.Item(139, 10).Value = 142876  # J139: 144220 -> 142876
$ws.Cells.Item(139, 11).Value = 0  # K139: 139000 -> 0
$ws.Cells.Item(139, 12).Value = 142876  # L139: 144220 -> 142876
$ws.Cells.Item(139, 13).Value = $null  # M139: -133860 -> None
$ws.Cells.Item(139, 14).Value = -153156  # N139: -154500 -> -153156
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 3034.4  # H5: 3328.4443 -> 3034.4
$ws.Cells.Item(5, 9).Value = 2873.25  # I5: 2873.5 -> 2873.25
$ws.Cells.Item(5, 10).Value = 3679  # J5: 6968 -> 3679
$ws.Cells.Item(5, 11).Value = 8619.75  # K5: 8620.5 -> 8619.75
$ws.Cells.Item(5, 12).Value = 11037  # L5: 20904 -> 11037
$ws.Cells.Item(5, 13).Value = -8507.75  # M5: -8508.5 -> -8507.75
$ws.Cells.Item(5, 14).Value = -11261  # N5: -21128 -> -11261
$ws.Cells.Item(75, 8).Value = 1333.5  # H75: 1665.6666 -> 1333.5
$ws.Cells.Item(75, 9).Value = 1442.6666  # I75: 1249.5 -> 1442.6666
$ws.Cells.Item(75, 10).Value = 1224.3334  # J75: 2498 -> 1224.3334
$ws.Cells.Item(75, 11).Value = 4327.9998  # K75: 3748.5 -> 4327.9998
$ws.Cells.Item(75, 12).Value = 3673.0002  # L75: 7494 -> 3673.0002
$ws.Cells.Item(75, 13).Value = -3329.9998  # M75: -2750.5 -> -3329.9998
$ws.Cells.Item(75, 14).Value = -5669.0002  # N75: -9490 -> -5669.0002
$ws.Cells.Item(78, 8).Value = 1333.5  # H78: 1665.6666 -> 1333.5
$ws.Cells.Item(78, 9).Value = 1442.6666  # I78: 1249.5 -> 1442.6666
$ws.Cells.Item(78, 10).Value = 1224.3334  # J78: 2498 -> 1224.3334
$ws.Cells.Item(78, 11).Value = 12983.9994  # K78: 11245.5 -> 12983.9994
$ws.Cells.Item(78, 12).Value = 11019.0006  # L78: 22482 -> 11019.0006
$ws.Cells.Item(78, 13).Value = -7991.999400000001  # M78: -6253.5 -> -7991.999400000001
$ws.Cells.Item(78, 14).Value = -21003.0006  # N78: -32466 -> -21003.0006
$ws.Cells.Item(132, 8).Value = 1423.6666  # H132: 1393.1875 -> 1423.6666
$ws.Cells.Item(132, 9).Value = 1045.5  # I132: 1067.8889 -> 1045.5
$ws.Cells.Item(132, 10).Value = 2180  # J132: 1811.4286 -> 2180
$ws.Cells.Item(132, 11).Value = 9409.5  # K132: 9611.000099999999 -> 9409.5
$ws.Cells.Item(132, 12).Value = 19620  # L132: 16302.8574 -> 19620
$ws.Cells.Item(132, 13).Value = -6879.5  # M132: -7081.000099999999 -> -6879.5
$ws.Cells.Item(132, 14).Value = -24680  # N132: -21362.8574 -> -24680
$ws.Cells.Item(135, 8).Value = 3034.4  # H135: 3328.4443 -> 3034.4
$ws.Cells.Item(135, 9).Value = 2873.25  # I135: 2873.5 -> 2873.25
$ws.Cells.Item(135, 10).Value = 3679  # J135: 6968 -> 3679
$ws.Cells.Item(135, 11).Value = 25859.25  # K135: 25861.5 -> 25859.25
$ws.Cells.Item(135, 12).Value = 33111  # L135: 62712 -> 33111
$ws.Cells.Item(135, 13).Value = -23324.25  # M135: -23326.5 -> -23324.25
$ws.Cells.Item(135, 14).Value = -38181  # N135: -67782 -> -38181
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(58, 8).Value = 19999  # H58: 0 -> 19999
$ws.Cells.Item(58, 10).Value = 19999  # J58: 0 -> 19999
$ws.Cells.Item(58, 12).Value = 19999  # L58: 0 -> 19999
$ws.Cells.Item(58, 14).Value = -20553  # N58: None -> -20553
$ws.Cells.Item(102, 8).Value = 4821.8887  # H102: 5174.75 -> 4821.8887
$ws.Cells.Item(102, 9).Value = 3699.5715  # I102: 3983 -> 3699.5715
$ws.Cells.Item(102, 11).Value = 3699.5715  # K102: 3983 -> 3699.5715
$ws.Cells.Item(102, 13).Value = -2077.5715  # M102: -2361 -> -2077.5715
$ws.Cells.Item(126, 8).Value = 2309.28  # H126: 2524.111 -> 2309.28
$ws.Cells.Item(126, 9).Value = 1695.4667  # I126: 1766.7273 -> 1695.4667
$ws.Cells.Item(126, 10).Value = 3230  # J126: 3714.2856 -> 3230
$ws.Cells.Item(126, 11).Value = 5086.4001  # K126: 5300.1819 -> 5086.4001
$ws.Cells.Item(126, 12).Value = 9690  # L126: 11142.8568 -> 9690
$ws.Cells.Item(126, 13).Value = -2616.4001  # M126: -2830.1819 -> -2616.4001
$ws.Cells.Item(126, 14).Value = -14630  # N126: -16082.8568 -> -14630
$ws.Cells.Item(132, 8).Value = 2867.3572  # H132: 3065.1538 -> 2867.3572
$ws.Cells.Item(132, 9).Value = 2014.4  # I132: 2094.2222 -> 2014.4
$ws.Cells.Item(132, 10).Value = 4999.75  # J132: 5249.75 -> 4999.75
$ws.Cells.Item(132, 11).Value = 6043.200000000001  # K132: 6282.6666 -> 6043.200000000001
$ws.Cells.Item(132, 12).Value = 14999.25  # L132: 15749.25 -> 14999.25
$ws.Cells.Item(132, 13).Value = -3513.200000000001  # M132: -3752.6666 -> -3513.200000000001
$ws.Cells.Item(132, 14).Value = -20059.25  # N132: -20809.25 -> -20059.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6120.1  # H7: 6633 -> 6120.1
$ws.Cells.Item(7, 9).Value = 2546.25  # I7: 2695.1428 -> 2546.25
$ws.Cells.Item(7, 11).Value = 2546.25  # K7: 2695.1428 -> 2546.25
$ws.Cells.Item(7, 13).Value = -2434.25  # M7: -2583.1428 -> -2434.25
$ws.Cells.Item(55, 8).Value = 856.82355  # H55: 898.75 -> 856.82355
$ws.Cells.Item(55, 9).Value = 999.1111  # I55: 918.8 -> 999.1111
$ws.Cells.Item(55, 10).Value = 696.75  # J55: 865.3333 -> 696.75
$ws.Cells.Item(55, 11).Value = 999.1111  # K55: 918.8 -> 999.1111
$ws.Cells.Item(55, 12).Value = 696.75  # L55: 865.3333 -> 696.75
$ws.Cells.Item(55, 13).Value = -826.1111  # M55: -745.8 -> -826.1111
$ws.Cells.Item(55, 14).Value = -1042.75  # N55: -1211.3333 -> -1042.75
$ws.Cells.Item(82, 8).Value = 4093.9714  # H82: 4436.3125 -> 4093.9714
$ws.Cells.Item(82, 9).Value = 4758.1333  # I82: 4398.2354 -> 4758.1333
$ws.Cells.Item(82, 10).Value = 3595.85  # J82: 4479.467 -> 3595.85
$ws.Cells.Item(82, 11).Value = 4758.1333  # K82: 4398.2354 -> 4758.1333
$ws.Cells.Item(82, 12).Value = 3595.85  # L82: 4479.467 -> 3595.85
$ws.Cells.Item(82, 13).Value = -4397.1333  # M82: -4037.2354 -> -4397.1333
$ws.Cells.Item(82, 14).Value = -4317.85  # N82: -5201.467 -> -4317.85
$ws.Cells.Item(85, 8).Value = 4093.9714  # H85: 4436.3125 -> 4093.9714
$ws.Cells.Item(85, 9).Value = 4758.1333  # I85: 4398.2354 -> 4758.1333
$ws.Cells.Item(85, 10).Value = 3595.85  # J85: 4479.467 -> 3595.85
$ws.Cells.Item(85, 11).Value = 4758.1333  # K85: 4398.2354 -> 4758.1333
$ws.Cells.Item(85, 12).Value = 3595.85  # L85: 4479.467 -> 3595.85
$ws.Cells.Item(85, 13).Value = -3510.1333  # M85: -3150.2354 -> -3510.1333
$ws.Cells.Item(85, 14).Value = -6091.85  # N85: -6975.467 -> -6091.85
$ws.Cells.Item(122, 8).Value = 3779.375  # H122: 3962.1428 -> 3779.375
$ws.Cells.Item(122, 9).Value = 3260  # I122: 3640 -> 3260
$ws.Cells.Item(122, 11).Value = 9780  # K122: 10920 -> 9780
$ws.Cells.Item(122, 13).Value = -7330  # M122: -8470 -> -7330
$ws.Cells.Item(126, 8).Value = 6120.1  # H126: 6633 -> 6120.1
$ws.Cells.Item(126, 9).Value = 2546.25  # I126: 2695.1428 -> 2546.25
$ws.Cells.Item(126, 11).Value = 7638.75  # K126: 8085.428400000001 -> 7638.75
$ws.Cells.Item(126, 13).Value = -5168.75  # M126: -5615.428400000001 -> -5168.75
$ws.Cells.Item(132, 8).Value = 7654.364  # H132: 7266.5 -> 7654.364
$ws.Cells.Item(132, 9).Value = 7654.364  # I132: 7266.5 -> 7654.364
$ws.Cells.Item(132, 11).Value = 22963.092  # K132: 21799.5 -> 22963.092
$ws.Cells.Item(132, 13).Value = -20433.092  # M132: -19269.5 -> -20433.092
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 25995  # H51: 25992.5 -> 25995
$ws.Cells.Item(51, 10).Value = 25995  # J51: 25992.5 -> 25995
$ws.Cells.Item(51, 12).Value = 25995  # L51: 25992.5 -> 25995
$ws.Cells.Item(51, 14).Value = -27015  # N51: -27012.5 -> -27015
$ws.Cells.Item(70, 8).Value = 29995  # H70: 33323.332 -> 29995
$ws.Cells.Item(70, 10).Value = 39990  # J70: 39985 -> 39990
$ws.Cells.Item(70, 12).Value = 39990  # L70: 39985 -> 39990
$ws.Cells.Item(70, 14).Value = -40620  # N70: -40615 -> -40620
$ws.Cells.Item(73, 8).Value = 29995  # H73: 33323.332 -> 29995
$ws.Cells.Item(73, 10).Value = 39990  # J73: 39985 -> 39990
$ws.Cells.Item(73, 12).Value = 39990  # L73: 39985 -> 39990
$ws.Cells.Item(73, 14).Value = -42174  # N73: -42169 -> -42174
$ws.Cells.Item(122, 8).Value = 3617.5417  # H122: 3692.2173 -> 3617.5417
$ws.Cells.Item(122, 9).Value = 2353.111  # I122: 2409.75 -> 2353.111
$ws.Cells.Item(122, 11).Value = 7059.333  # K122: 7229.25 -> 7059.333
$ws.Cells.Item(122, 13).Value = -4609.333  # M122: -4779.25 -> -4609.333
